$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily data point (12.01.2021) was reported. It belongs above the
# existing 11.01.2021 row (row 234), so insert a fresh row there and push
# everything from row 234 down (234->235, ..., 239->240).
$ws.Rows.Item(234).Insert()

# Column A holds the date as text (matches the rest of the column, which is
# stored as text rather than a real Excel date), so force text formatting
# before writing the value to stop it being auto-converted to a date serial.
$ws.Cells.Item(234, 1).NumberFormat = "@"
$ws.Cells.Item(234, 1).Value = "12.01.2021"
$ws.Cells.Item(234, 2).Value = 116200
$ws.Cells.Item(234, 3).Value = 548818
$ws.Cells.Item(234, 4).Value = 4373
$ws.Cells.Item(234, 5).Value = 82822
$ws.Cells.Item(234, 6).Value = 29005
$ws.Cells.Item(234, 7).Value = 0
